# Add the new "Our study guides..." intro paragraph (style FirstParagraph)
# right after the Title paragraph, and add a "contents" bookmark just after
# it (mirrors the bookmarkStart/bookmarkEnd pair that sits right before
# <w:sectPr> in the target document).

$d = $word.ActiveDocument

# --- 1. Insert the new paragraph using a literal OOXML fragment so it comes
#        out clean (no stray w:rsidP/w:rsidRPr attributes, and the run text
#        keeps xml:space="preserve" exactly like the target). -------------
$insertionPoint = $d.Content
$insertionPoint.Collapse(0)

$paraXml = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr>' +
              '<w:r><w:t xml:space="preserve">Our study guides, written for students by students on a specific area of mathematics or statistics, can be found here.</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$insertionPoint.InsertXML($paraXml)

# --- 2. Add the "contents" bookmark right after the new paragraph's text.
#        Bookmarks.Add mis-handles a collapsed range sitting exactly on a
#        paragraph-mark boundary (it "resets" to the top of the document),
#        so we dodge that edge case: type a throw-away character at the
#        very end of the document, drop the bookmark immediately before
#        that character (a perfectly ordinary mid-run position), then
#        delete the throw-away character again. The bookmark stays put,
#        now sitting correctly at the end of the paragraph we just added. -
$d.Content.InsertAfter("X")
$endPos = $d.Content.End
$bookmarkPos = $endPos - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("contents", $bookmarkRange)
$d.Range($bookmarkPos, $bookmarkPos + 1).Delete()

Write-Output "paragraphs=$($d.Paragraphs.Count) bookmarks=$($d.Bookmarks.Count)"
